# Apply resume content edits via Word COM-interop Find/Replace.
$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $range = $d.Content
    $ok = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

# 1. Summary paragraph — consolidate into one sentence-pair.
$old1 = "Experienced engineering leader with over 20 years of experience in software development, technical strategy & consulting.  As a hands-on technical strategist, software architect, engineer, blogger, consultant & OSS contributor, I have used my extensive leadership & technical experience to implement innovative solutions across several industries."
$new1 = "Versatile engineering leader with deep expertise in software development, technical strategy, and consulting. A hands-on architect, engineer, blogger, OSS contributor, and trusted advisor, I bring a proven track record of delivering innovative, scalable solutions across diverse industries."
Replace-Text $old1 $new1

# 2. Production Fundamentals OpenHack bullet — drop the spell-check markup around "OpenHack".
$old2 = "Part of a team of engineers that created the Production Fundamentals OpenHack (OH).  The OH provided hands"
$new2 = "Part of a team of engineers that created the Production Fundamentals OpenHack (OH).  The OH provided hands"
Replace-Text $old2 $new2

# 3. Executive Dashboard bullet — shortened wording.
$old3 = "Responsible for development, architecture, release management & the overall technical roadmap for the Executive Dashboard Application."
$new3 = "Directed architecture, development, and release strategy for an enterprise Executive Dashboard application."
Replace-Text $old3 $new3

# 4. Healthcare claims platform bullet — shortened wording.
$old4 = "Helped client move into the healthcare industry by performing work on their medical claims platform."
$new4 = "Enabled client’s entry into healthcare by enhancing medical claims platform and designing audit database."
Replace-Text $old4 $new4

# 5. Merge "I worked on..." bullet with new text, then delete the following
#    "Helped lead a multi-million-dollar..." bullet paragraph entirely.
$old5 = "I worked on the customer’s medical claims file parsing & data translation enterprise application.  Also helped design medical claims database used for the medical claims audit portion of the software."
$new5 = "Delivered mission-critical, multi-million-dollar state tax collection system; led development and mentored team members."
Replace-Text $old5 $new5

# Now remove the obsolete paragraph that followed (its content is now redundant).
# Deleting the paragraph's own Range (which includes its trailing paragraph
# mark) removes the whole <w:p> instead of leaving an empty bullet behind.
$deleteMarker = "Helped lead a multi-million-dollar"
$removed = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$deleteMarker*") {
        $p.Range.Delete() | Out-Null
        $removed = $true
        break
    }
}
if (-not $removed) {
    throw "Could not find paragraph to delete"
}
